$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update as explicit text, preserving the original
# "inlineStr" / text semantics (these columns hold formatted price/
# percentage strings, not numeric values) and the cells default style.
$updates = @(
    @{ Cell = 'D2'; Value = '63.108.51' },
    @{ Cell = 'E2'; Value = '  -1.46%  ' },
    @{ Cell = 'D3'; Value = '3.057.01' },
    @{ Cell = 'E3'; Value = '  -3.01%  ' },
    @{ Cell = 'E4'; Value = '  -0.26%  ' },
    @{ Cell = 'D5'; Value = '589.80' },
    @{ Cell = 'D6'; Value = '153.03' },
    @{ Cell = 'E6'; Value = '  +4.67%  ' },
    @{ Cell = 'E7'; Value = '  -0.14%  ' },
    @{ Cell = 'E8'; Value = '  +2.95%  ' },
    @{ Cell = 'D9'; Value = '3.062.68' },
    @{ Cell = 'E9'; Value = '  -2.51%  ' },
    @{ Cell = 'E10'; Value = '  -3.52%  ' },
    @{ Cell = 'E11'; Value = '  -0.95%  ' },
    @{ Cell = 'E12'; Value = '  -0.08%  ' },
    @{ Cell = 'E13'; Value = '  -2.79%  ' },
    @{ Cell = 'D14'; Value = '37.12' },
    @{ Cell = 'E14'; Value = '  -0.27%  ' },
    @{ Cell = 'E15'; Value = '  -1.95%  ' },
    @{ Cell = 'D16'; Value = '3.563.66' },
    @{ Cell = 'E16'; Value = '  -3.01%  ' },
    @{ Cell = 'E17'; Value = '  -1.17%  ' },
    @{ Cell = 'D18'; Value = '63.173.51' },
    @{ Cell = 'E18'; Value = '  -1.14%  ' },
    @{ Cell = 'D19'; Value = '3.061.20' },
    @{ Cell = 'E19'; Value = '  -2.75%  ' },
    @{ Cell = 'D20'; Value = '474.19' },
    @{ Cell = 'E20'; Value = '  +1.44%  ' },
    @{ Cell = 'D21'; Value = '14.62' },
    @{ Cell = 'E21'; Value = '  +1.71%  ' },
    @{ Cell = 'E22'; Value = '  -2.21%  ' },
    @{ Cell = 'D23'; Value = '7.53' },
    @{ Cell = 'E23'; Value = '  +0.24%  ' },
    @{ Cell = 'D24'; Value = '2.38' },
    @{ Cell = 'E24'; Value = '  +2.12%  ' },
    @{ Cell = 'E25'; Value = '  -0.41%  ' },
    @{ Cell = 'D26'; Value = '81.04' },
    @{ Cell = 'E26'; Value = '  -0.23%  ' },
    @{ Cell = 'E27'; Value = '  -0.28%  ' },
    @{ Cell = 'D28'; Value = '9.98' },
    @{ Cell = 'E28'; Value = '  +2.39%  ' },
    @{ Cell = 'E29'; Value = '  -1.35%  ' },
    @{ Cell = 'D30'; Value = '0.999' },
    @{ Cell = 'E30'; Value = '  -0.26%  ' },
    @{ Cell = 'E31'; Value = '  -1.91%  ' },
    @{ Cell = 'E32'; Value = '  -2.07%  ' },
    @{ Cell = 'D34'; Value = '27.18' },
    @{ Cell = 'E34'; Value = '  -1.87%  ' },
    @{ Cell = 'D35'; Value = '0.0₃0840' },
    @{ Cell = 'E35'; Value = '  +0.24%  ' },
    @{ Cell = 'D36'; Value = '1.05' },
    @{ Cell = 'E36'; Value = '  -2.09%  ' },
    @{ Cell = 'D37'; Value = '6.10' },
    @{ Cell = 'E37'; Value = '  -1.05%  ' },
    @{ Cell = 'E38'; Value = '  +1.94%  ' },
    @{ Cell = 'D39'; Value = '2.21' },
    @{ Cell = 'E39'; Value = '  -4.92%  ' },
    @{ Cell = 'D40'; Value = '9.28' },
    @{ Cell = 'E40'; Value = '  +0.77%  ' },
    @{ Cell = 'D41'; Value = '50.37' },
    @{ Cell = 'E41'; Value = '  -2.05%  ' },
    @{ Cell = 'D42'; Value = '443.57' },
    @{ Cell = 'E42'; Value = '  -4.39%  ' },
    @{ Cell = 'E43'; Value = '  -3.11%  ' },
    @{ Cell = 'B44'; Value = 'VeChain' },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' },
    @{ Cell = 'D44'; Value = '0.0362' },
    @{ Cell = 'E44'; Value = '  -2.64%  ' },
    @{ Cell = 'B45'; Value = 'Arweave' },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar' },
    @{ Cell = 'D45'; Value = '40.09' },
    @{ Cell = 'E45'; Value = '  -0.38%  ' },
    @{ Cell = 'E46'; Value = '  +2.11%  ' },
    @{ Cell = 'D47'; Value = '2.792.77' },
    @{ Cell = 'E47'; Value = '  -4.56%  ' },
    @{ Cell = 'D48'; Value = '131.37' },
    @{ Cell = 'E48'; Value = '  +1.81%  ' },
    @{ Cell = 'D50'; Value = '25.10' },
    @{ Cell = 'E50'; Value = '  +3.55%  ' },
    @{ Cell = 'D51'; Value = '2.26' },
    @{ Cell = 'E51'; Value = '  +0.21%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}
